$d = $word.ActiveDocument

# The original single paragraph reads:
#   "This is about venus" + bookmarkStart/bookmarkEnd ("_GoBack")
#
# Target layout (per the diff):
#   1) "This is about venus"
#   2) <empty paragraph>
#   3) "More abouts venus"
#   4) <empty paragraph containing just the _GoBack bookmark>

# Step 1: split the lone paragraph so the trailing bookmark becomes its own,
# otherwise-empty paragraph at the end of the document. This preserves the
# bookmark exactly where it already is while detaching it from the text run.
$d.Content.Find.Execute("venus", $true, $false, $false, $false, $false, $true, 1, $false, "venus^p", 2)

# Step 2: insert the new empty paragraph and the "More abouts venus"
# paragraph right after the original text, and before the now-isolated
# bookmark paragraph. Using InsertXML lets us add a truly empty <w:p/>
# rather than one carrying a placeholder run.
$r = $d.Range(19, 19)
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newXml = "<w:p $xmlNs/><w:p $xmlNs><w:r><w:t>More abouts venus</w:t></w:r></w:p>"
$r.InsertXML($newXml)
